$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.889.84'
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('D3').Value = '2.243.92'
$ws.Range('E3').Value = '  -1.87%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '115.54'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '301.67'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +13.52%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.631'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.90%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '46.41'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.96%  '
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('B12').Value = 'OKB'
$ws.Range('C12').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '56.87'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.18%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '9.10'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.00%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.105'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.65%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.39'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.35%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.892'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('B17').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C17').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D17').Value = '2.581.87'
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.266.89'
$ws.Range('E18').Value = '  -0.92%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '42.710.49'
$ws.Range('E19').Value = '  -2.07%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.55'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +11.41%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0000108'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '74.21'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.39%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.54'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +23.14%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.36'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.63%  '
$ws.Range('B25').Value = 'BitcoinCash'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '232.80'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.58%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.41'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.22'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +5.52%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.94%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '40.16'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.09%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.23'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.77%  '
$ws.Range('B31').Value = 'WEMIXToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.27'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.57%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '175.68'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.10%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '21.29'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -2.07%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0909'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.58'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +16.27%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.64'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.128'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.46%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.76'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.48%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0374'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.106'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.61'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.42%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '72.51'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.49%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '13.53'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -5.80%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.238'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.04%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.35'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.85%  '
$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '5.60'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -6.45%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.34'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.75%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '107.20'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +6.70%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.62'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0990'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.85%  '
